$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 60.71984751580217
$ws.Range("B3").Value = 66.07015436451078
$ws.Range("B4").Value = 81.12193842253824
$ws.Range("H5").Value = 96.18691285694622
$ws.Range("H6").Value = 96.20100698668097
$ws.Range("H7").Value = 96.18082039348931
$ws.Range("C8").Value = 86.09332943401988
$ws.Range("C9").Value = 86.56770807175711
$ws.Range("C10").Value = 86.81184556348657
$ws.Range("D11").Value = 99.21958148675363
$ws.Range("D12").Value = 99.27802673591547
$ws.Range("D13").Value = 99.39171778791176
$ws.Range("E14").Value = 98.93658064173701
$ws.Range("E15").Value = 98.99412345848863
$ws.Range("E16").Value = 98.9458616453217
$ws.Range("F17").Value = 98.47513980388757
$ws.Range("F18").Value = 98.55337953282209
$ws.Range("F19").Value = 98.45312589083835
$ws.Range("G20").Value = 97.54987876951452
$ws.Range("G21").Value = 97.65553546247817
$ws.Range("G22").Value = 97.58766673988303
$ws.Range("B23").Value = 82.25554390024789
$ws.Range("B24").Value = 86.17557709739998
$ws.Range("H25").Value = 96.21438089214904
$ws.Range("H26").Value = 96.14469766926345
$ws.Range("C27").Value = 86.19615976744174
$ws.Range("C28").Value = 86.45685250944567
$ws.Range("D29").Value = 99.33277181411147
$ws.Range("D30").Value = 99.30790263815177
$ws.Range("E31").Value = 98.95350190457305
$ws.Range("E32").Value = 98.94748754800725
$ws.Range("F33").Value = 98.49284508385018
$ws.Range("F34").Value = 98.55915295571171
$ws.Range("G35").Value = 97.51807646100171
$ws.Range("G36").Value = 97.54644667639212
$ws.Range("B37").Value = 73.70854751749948
$ws.Range("B38").Value = 81.20994996761762
$ws.Range("H39").Value = 96.14283544976692
$ws.Range("H40").Value = 96.05325097265769
$ws.Range("C41").Value = 86.14915281209417
$ws.Range("C42").Value = 86.33137278515588
$ws.Range("D43").Value = 99.38644295733721
$ws.Range("D44").Value = 99.3298419525742
$ws.Range("E45").Value = 98.98196196311666
$ws.Range("E46").Value = 98.85748858113767
$ws.Range("F47").Value = 98.55592964630206
$ws.Range("F48").Value = 98.40224172209965
$ws.Range("G49").Value = 97.5259085391735
$ws.Range("G50").Value = 97.57467590858005
